$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet structure: rename the current "Container" sheet in place to
#    "Provenance" (it keeps its sheetId / r:id). A brand new "Container"
#    sheet is appended right after it -- duplicated from the original
#    Container sheet (before the rename) so it inherits the same
#    boilerplate (sheetPr/sheetFormatPr/pageMargins) as the rest of the
#    workbook instead of generic Worksheets.Add() defaults.
# ---------------------------------------------------------------------------
$provenance = $wb.Worksheets.Item("Container")
$provenance.Copy($null, $provenance)
$provenance.Name = "Provenance"

$container = $wb.Worksheets.Item("Container (2)")
$container.Name = "Container"

# ---------------------------------------------------------------------------
# 2. Food sheet: new header row (A1:Z1) + refreshed data validations.
# ---------------------------------------------------------------------------
$food = $wb.Worksheets.Item("Food")

$foodHeaders = @(
    "food_primary_type",
    "food_primary_type_label",
    "food_upc_code",
    "food_preservation_state",
    "food_storage_temperature_state",
    "food_ripeness_state",
    "food_cooking_method",
    "food_acquisition_city",
    "food_acquisition_country",
    "food_acquisition_country_subdivision",
    "food_acquisition_date",
    "food_acquisition_location_type",
    "food_acquisition_latitude",
    "food_acquisition_longitude",
    "food_acquisition_agent_name",
    "food_acquisition_organization",
    "food_distributor_city",
    "food_distributor_country",
    "food_distributor_country_subdivision",
    "food_expiration_date",
    "food_category_label",
    "food_additional_types",
    "id",
    "laboratory_sample_id",
    "laboratory_sample_aliquot_id",
    "laboratory_sample_batch_id"
)
for ($i = 0; $i -lt $foodHeaders.Length; $i++) {
    $food.Cells.Item(1, $i + 1).Value = $foodHeaders[$i]
}

# Drop the two old validations (old D + old I) before re-adding the new set.
$food.Range("D2:D1048576").Validation.Delete()
$food.Range("I2:I1048576").Validation.Delete()

$foodValidations = @{
    "D" = "air-dried,artificially dried,brined,candied,canned,cured,dried,fermented,freeze-dried,fresh,heat treated,irradiated,jellied,kippered,naturally dried,pasteurized,pickled,raw,shelf stable,sun-dried,ultraviolet light exposed";
    "E" = "chilled,foodsafe chilled,frozen,refrigerated";
    "F" = "ripe,overripe,unripe,slightly ripe";
    "L" = "field,fresh market,small grocery,supermarket,biobank,unknown,other"
}
foreach ($col in @("D", "E", "F", "L")) {
    $rng = $food.Range("$($col)2:$($col)1048576")
    $v = $rng.Validation
    $v.Add(3, 1, 1, '"' + $foodValidations[$col] + '"')
    $v.IgnoreBlank = $true
    $v.InCellDropdown = $true
    $v.ShowInput = $false
    $v.ShowError = $false
}

# ---------------------------------------------------------------------------
# 3. Component sheet: new header row (A1:V1) + refreshed data validation.
# ---------------------------------------------------------------------------
$component = $wb.Worksheets.Item("Component")

$componentHeaders = @(
    "component_type",
    "component_type_label",
    "component_recorded_value",
    "component_measurement_unit",
    "component_data_points_number",
    "component_record_date",
    "component_analysis_date",
    "component_comment",
    "component_derivation_type",
    "component_limit_of_quantitation",
    "laboratory_sample_aggregation_minimum_measured_compound_value",
    "laboratory_sample_aggregation_maximum_measured_compound_value",
    "laboratory_sample_aggregation_median_measured_compound_value",
    "laboratory_sample_aggregation_measured_compound_standard_deviation",
    "analytical_analysis_measurement_protocol_doi",
    "analytical_analysis_measurement_method",
    "laboratory_conducting_analytical_analysis",
    "component_quality_control_remeasurement",
    "id",
    "laboratory_sample_id",
    "laboratory_sample_aliquot_id",
    "laboratory_sample_batch_id"
)
for ($i = 0; $i -lt $componentHeaders.Length; $i++) {
    $component.Cells.Item(1, $i + 1).Value = $componentHeaders[$i]
}

# Old validation lived on Q; new one lives on P.
$component.Range("Q2:Q1048576").Validation.Delete()
$vc = $component.Range("P2:P1048576").Validation
$vc.Add(3, 1, 1, '"HPLC,GLC,GC,Nephelometry,Gravimetric,Fluorometric,Kjeldahl"')
$vc.IgnoreBlank = $true
$vc.InCellDropdown = $true
$vc.ShowInput = $false
$vc.ShowError = $false

# ---------------------------------------------------------------------------
# 4. Provenance sheet (previously "Container"): brand new header row.
# ---------------------------------------------------------------------------
$provenanceHeaders = @(
    "dataset_label",
    "mifc_version_tag",
    "contributor_orcid",
    "organization_name",
    "id",
    "laboratory_sample_id",
    "laboratory_sample_aliquot_id",
    "laboratory_sample_batch_id"
)
for ($i = 0; $i -lt $provenanceHeaders.Length; $i++) {
    $provenance.Cells.Item(1, $i + 1).Value = $provenanceHeaders[$i]
}

# ---------------------------------------------------------------------------
# 5. Container sheet (brand new): the three linked-table names.
# ---------------------------------------------------------------------------
$containerHeaders = @("foods", "components", "provenances")
for ($i = 0; $i -lt $containerHeaders.Length; $i++) {
    $container.Cells.Item(1, $i + 1).Value = $containerHeaders[$i]
}
